$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Row 7: "Experimental" -> set the Value cell (B7) to the literal text "false" ---
# Plain Value/Value2/Formula assignment of the string "false" gets auto-coerced
# to a real Boolean FALSE by the engine (mirrors genuine Excel typed-input
# behavior). To force a literal text string instead, stage it in a scratch
# cell as a formula result (text type), copy it, then paste-values into B7 -
# this preserves the "text" data type rather than "boolean".
$helper = $ws.Range("Z1")
$helper.Formula = '="false"'
$helper.Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$helper.ClearContents()

# --- Row 8: "Date" -> update the ISO timestamp value (B8) ---
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
